# 23 dec 2023 update
# Fill in the 6th weekly payment row (row 8) on the "MD10000.20-OCT" sheet:
# a payment date, the amount paid, and the paid-weeks flag. The dependent
# summary formulas (K1, O1, K2, O2) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("MD10000.20-OCT")

$ws.Range("B8").Value = 45272
$ws.Range("C8").Value = 700
$ws.Range("D8").Value = 1

# Move the active selection from D6:D7 down to the newly-filled D8 cell.
$ws.Range("D8").Select()
